$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number for each data row (2..171).
# This update bumps that date by one day (45178 -> 45179) for every row.
for ($r = 2; $r -le 171; $r++) {
    $ws.Cells.Item($r, 3).Value = 45179
}
